# unitTest_repeatUntil.xlsx - add new "outputToCloud(resource)" function under the
# existing "base" category, and add a brand new "text" category that holds the
# new "spellCheck(var,profile,text)" function.
#
# The '#system' sheet (sheet1) lays out one category per column: row 1 holds the
# category name (which doubles as the column header used by the data-validation
# dropdowns on the visible sheet), and the cells below list every function that
# belongs to that category, in alphabetical order. Each category also has a
# matching workbook-level defined Name, e.g. base => '#system'!$E$2:$E$38.
#
# This change:
#   1. Inserts "outputToCloud(resource)" into the "base" column (E), in its
#      alphabetically-sorted position (just before "prependText").
#   2. Inserts a new category "text" between "step" and "web" in the category
#      list (column A), in alphabetical order.
#   3. Makes room for the new "text" category by shifting the "web", "webalert",
#      "webcookie", "ws", "ws.async" and "xml" columns one column to the right
#      (Y..AD -> Z..AE), then fills the now-empty column Y with the "text"
#      category header and its single function "spellCheck(var,profile,text)".
#   4. Updates the defined Names so they keep pointing at the right ranges, and
#      adds a new Name "text" for the new category.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------------
# 1) "base" column (E): insert outputToCloud(resource) before prependText(...)
#    E2:E38 (37 entries) -> E2:E39 (38 entries)
# ---------------------------------------------------------------------------
$baseShiftSrc = $ws.Range("E22:E38")
$baseShiftDst = $ws.Range("E23:E39")
$baseShiftDst.Value = $baseShiftSrc.Value()
$ws.Range("E22").Value = "outputToCloud(resource)"

# ---------------------------------------------------------------------------
# 2) "target" column (A): insert the new "text" category name before "web"
#    A2:A30 (29 entries) -> A2:A31 (30 entries)
# ---------------------------------------------------------------------------
$targetShiftSrc = $ws.Range("A25:A30")
$targetShiftDst = $ws.Range("A26:A31")
$targetShiftDst.Value = $targetShiftSrc.Value()
$ws.Range("A25").Value = "text"

# ---------------------------------------------------------------------------
# 3) Shift the web / webalert / webcookie / ws / ws.async / xml columns one
#    column to the right (Y..AD -> Z..AE) to make room for the "text" column,
#    then populate the freed-up column Y with the new "text" category.
# ---------------------------------------------------------------------------
$colShiftSrc = $ws.Range("Y1:AD129")
$colShiftDst = $ws.Range("Z1:AE129")
$colShiftDst.Value = $colShiftSrc.Value()

$ws.Range("Y1").Value = "text"
$ws.Range("Y2").Value = "spellCheck(var,profile,text)"
$ws.Range("Y3:Y129").ClearContents()

# ---------------------------------------------------------------------------
# 4) Update the defined Names to reflect the new ranges, and add the new
#    "text" Name.
# ---------------------------------------------------------------------------
$wb.Names.Item("base").RefersTo = "='#system'!`$E`$2:`$E`$39"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$31"
$wb.Names.Item("web").RefersTo = "='#system'!`$Z`$2:`$Z`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AC`$2:`$AC`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AD`$2:`$AD`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AE`$2:`$AE`$27"

$wb.Names.Add("text", "='#system'!`$Y`$2:`$Y`$2")
